$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 2, shifting existing data (rows 2-8) down to rows 4-10.
$ws.Range("A2:K3").EntireRow.Insert()

# B2:C3 hold numeric-looking identifiers (account/CIN numbers) that must be
# kept as text (matching the rest of the column). Temporarily mark them as
# Text ("@") so the numeric-looking strings aren't coerced to numbers, then
# restore the default "Normal" style so no visible formatting change remains.
$ws.Range("B2:C3").NumberFormat = "@"

# Row 2: new entry - STE LOCATION
$ws.Range("A2").Value = "STE LOCATION "
$ws.Range("B2").Value = "31451"
$ws.Range("C2").Value = "313156456461638489461313"
$ws.Range("D2").Value = "MAARIF"
$ws.Range("E2").Value = "BP"
$ws.Range("F2").Value = "Point de vente"
$ws.Range("G2").Value = "512/CASA 2"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 5000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 5000

# Row 3: new entry - STE MAISON
$ws.Range("A3").Value = "STE MAISON "
$ws.Range("B3").Value = "56987"
$ws.Range("C3").Value = "322656131365484946461313"
$ws.Range("D3").Value = "2 MARS "
$ws.Range("E3").Value = "BMCE"
$ws.Range("F3").Value = "Point de vente"
$ws.Range("G3").Value = "512/CASA 2"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 5000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5000

# Restore default styling now that the values are safely stored as text.
$ws.Range("B2:C3").Style = "Normal"

# Row 4 (previously row 2 - STT22): fix G4 value
$ws.Range("G4").Value = "903/CASA ANFA/AV"

# Row 9 (previously row 7 - JAJA GAGA): fix J9 value
$ws.Range("J9").Value = 2800

# Row 10 (previously row 8 - MANAL LALA): fix J10 value
$ws.Range("J10").Value = 22500
